# "Generate Report for Handoff": a fresh handoff package was produced for the
# tracked file, so its generated identifier/hash and the handoff timestamps
# move forward, and the (now-superseded) handback info is reset.

$wb = $excel.ActiveWorkbook

$oldUuid = "d93a9d40-8581-46d8-8641-19ac82c3f509"
$newUuid = "a5c1953a-cb1d-4a8f-9183-a733b2dce419"
$oldHash = "c6ce4e430e783b9f8d29f3cb297000026aa2ad86"
$newHash = "54d2422e7736608e7930b91a8a704b8f2211fddb"
$resetDate = "0001-01-01 00:00:00"

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $ov.Range("A2").Value2.Replace($oldUuid, $newUuid)
$ov.Range("B2").Value = $ov.Range("B2").Value2.Replace($oldUuid, $newUuid)
$ov.Range("G2").Value = "2016-08-27 16:58:18"

# Update the hyperlink display text on B2 (target address is unchanged).
$ovLinkAddr = ""
foreach ($hl in $ov.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$2") { $ovLinkAddr = $hl.Address }
}
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $ovLinkAddr, "", "", "e2e\" + $newUuid + ".md")

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $zh.Range("A2").Value2.Replace($oldUuid, $newUuid)
$zh.Range("G2").Value = $zh.Range("G2").Value2.Replace($oldUuid, $newUuid).Replace($oldHash, $newHash)
$zh.Range("H2").Value = "2016-08-27 16:58:13"

# Keep only the A2 hyperlink (re-added with the refreshed display text);
# the I2 "latest target file" link goes away along with its cell content.
$zhLinkAddr = ""
foreach ($hl in $zh.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") { $zhLinkAddr = $hl.Address }
}
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhLinkAddr, "", "", $newUuid + ".md")

$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = $resetDate

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $de.Range("A2").Value2.Replace($oldUuid, $newUuid)
$de.Range("G2").Value = $de.Range("G2").Value2.Replace($oldUuid, $newUuid).Replace($oldHash, $newHash)
$de.Range("H2").Value = "2016-08-27 16:58:18"

$deLinkAddr = ""
foreach ($hl in $de.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") { $deLinkAddr = $hl.Address }
}
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deLinkAddr, "", "", $newUuid + ".md")

$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("K2").Value = $resetDate

# Column widths on the two locale sheets shrink now that the Latest
# Target File / Latest Handback File columns are empty (values chosen so
# the engine's character-width quantization lands on the same stored
# width Excel's AutoFit produced: ~18.65 / ~21.71 chars).
foreach ($ws in @($zh, $de)) {
    $ws.Columns.Item(9).ColumnWidth = 17.83
    $ws.Columns.Item(10).ColumnWidth = 20.83
}
